$wb = $excel.ActiveWorkbook

# --- Characteristics sheet: add a "Recovered" population-size row so the
#     recovered compartment isn't left with negative people. ---
$ws = $wb.Worksheets.Item("Characteristics")

# Copy formatting from the row above (row 10) into the new row 11, in two
# pieces so we don't pick up column F's default style on a cell that has
# no counterpart in row 10.
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("G10:M10").Copy()
$ws.Range("G11:M11").PasteSpecial(-4122)

$ws.Range("A11").Value = "rec_label"
$ws.Range("B11").Value = "Recovered"
$ws.Range("C11").Value = "sh_cases"
$ws.Range("E11").Value = 5000
$ws.Range("G11").Value = "rec"
$ws.Range("I11").Value = "rec"

$ws.Range("A11").Select()

# --- Update which sheet/tab is active: Characteristics becomes the
#     selected tab, and "Databook Sheet Names" (previously selected) no
#     longer is. ---
$wb.Worksheets.Item("Databook Sheet Names").Activate()
$ws.Activate()
